# Applies numeric value updates to the Kraken_Profits price-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as produced by the scheduled price-refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2833.3333
$ws.Range("J69").Value = 2999.5
$ws.Range("L69").Value = 8998.5
$ws.Range("N69").Value = -10746.5
$ws.Range("H72").Value = 2833.3333
$ws.Range("J72").Value = 2999.5
$ws.Range("L72").Value = 26995.5
$ws.Range("N72").Value = -35731.5
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
$ws.Range("H130").Value = 96383.60000000001
$ws.Range("J130").Value = 96383.60000000001
$ws.Range("L130").Value = 96383.60000000001
$ws.Range("N130").Value = -106423.6
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6750
$ws.Range("J45").Value = 7500
$ws.Range("L45").Value = 7500
$ws.Range("N45").Value = -8254
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H63").Value = 1749.5
$ws.Range("I63").Value = 1749.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1749.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1063.5
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 1749.5
$ws.Range("I66").Value = 1749.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8747.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5315.5
$ws.Range("N66").Value = $null
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960
$ws.Range("H129").Value = 39999.855
$ws.Range("J129").Value = 39999.855
$ws.Range("L129").Value = 39999.855
$ws.Range("N129").Value = -49999.855
$ws.Range("H133").Value = 99987
$ws.Range("J133").Value = 99987
$ws.Range("L133").Value = 99987
$ws.Range("N133").Value = -105047
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H107").Value = 2135.7778
$ws.Range("I107").Value = 1902.75
$ws.Range("K107").Value = 1902.75
$ws.Range("M107").Value = 17.25
$ws.Range("H115").Value = 90000
$ws.Range("I115").Value = 90000
$ws.Range("K115").Value = 90000
$ws.Range("M115").Value = -88433
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H129").Value = 75000
$ws.Range("J129").Value = 75000
$ws.Range("L129").Value = 75000
$ws.Range("N129").Value = -85000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6436.5835
$ws.Range("J31").Value = 7125.7
$ws.Range("L31").Value = 7125.7
$ws.Range("N31").Value = -7715.7
$ws.Range("H34").Value = 6436.5835
$ws.Range("J34").Value = 7125.7
$ws.Range("L34").Value = 7125.7
$ws.Range("N34").Value = -7529.7
$ws.Range("H53").Value = 44500
$ws.Range("J53").Value = 44500
$ws.Range("L53").Value = 44500
$ws.Range("N53").Value = -45714
$ws.Range("H98").Value = 95000
$ws.Range("I98").Value = 95000
$ws.Range("K98").Value = 95000
$ws.Range("M98").Value = -92754
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 374.75
$ws.Range("I92").Value = 799
$ws.Range("J92").Value = 233.33333
$ws.Range("K92").Value = 2397
$ws.Range("L92").Value = 699.99999
$ws.Range("M92").Value = -1149
$ws.Range("N92").Value = -3195.99999
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 99973
$ws.Range("J110").Value = 99973
$ws.Range("L110").Value = 99973
$ws.Range("N110").Value = -108153
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 8004
$ws.Range("I3").Value = 8004
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 8004
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -7892
$ws.Range("N3").Value = $null
$ws.Range("H14").Value = 7004
$ws.Range("I14").Value = 7004
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 7004
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -6832
$ws.Range("N14").Value = $null
$ws.Range("H15").Value = 8004
$ws.Range("I15").Value = 8004
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 8004
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -7834
$ws.Range("N15").Value = $null
$ws.Range("H16").Value = 2528.2
$ws.Range("J16").Value = 2139.6667
$ws.Range("L16").Value = 2139.6667
$ws.Range("N16").Value = -2479.6667
$ws.Range("H82").Value = 2325
$ws.Range("I82").Value = 2190
$ws.Range("K82").Value = 2190
$ws.Range("M82").Value = -1829
$ws.Range("H85").Value = 2325
$ws.Range("I85").Value = 2190
$ws.Range("K85").Value = 2190
$ws.Range("M85").Value = -942
$ws.Range("H100").Value = 5893.143
$ws.Range("I100").Value = 5893.143
$ws.Range("K100").Value = 5893.143
$ws.Range("M100").Value = -5352.143
$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815
$ws.Range("H125").Value = 39000
$ws.Range("J125").Value = 39000
$ws.Range("L125").Value = 39000
$ws.Range("N125").Value = -48840
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 30611.6
$ws.Range("I81").Value = 30611.6
$ws.Range("K81").Value = 61223.2
$ws.Range("M81").Value = -60162.2
$ws.Range("H84").Value = 30611.6
$ws.Range("I84").Value = 30611.6
$ws.Range("K84").Value = 306116
$ws.Range("M84").Value = -300812
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
$ws.Range("H128").Value = 60857.5
$ws.Range("J128").Value = 60857.5
$ws.Range("L128").Value = 60857.5
$ws.Range("N128").Value = -70817.5
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
